$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: make the "6-4 explore more about $group & $project" heading
# paragraph bold (paragraph mark + run).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$headingIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "6-4 explore more about*") {
        $headingIndex = $i
        break
    }
}
if ($headingIndex -eq -1) {
    throw "Could not find the '6-4 explore more about...' paragraph"
}
$d.Paragraphs.Item($headingIndex).Range.Bold = 1

# ------------------------------------------------------------------
# Edit 2: insert the new "6-5 Explore $group with $unwind aggregation
# stage" section right after the final "    }" / before the final "])"
# of the 6-4 section (i.e. right after the paragraph that follows the
# one containing "rangeBetweenMaxAndMin").
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*rangeBetweenMaxAndMin*") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not find the 'rangeBetweenMaxAndMin' paragraph"
}
# anchorIndex      -> "            rangeBetweenMaxAndMin: {...}"
# anchorIndex + 1  -> "        }"
# anchorIndex + 2  -> "    }"
# anchorIndex + 3  -> "])"      <- insertion point is right before this one
#
# NOTE: collapsing the *previous* paragraph's range to its end and
# inserting there causes the engine to merge the new content into that
# paragraph (losing its original text). Collapsing the *following*
# paragraph's range to its start instead inserts cleanly between the
# two paragraphs, preserving both.
$insertBefore = $d.Paragraphs.Item($anchorIndex + 3)
if ($insertBefore.Range.Text -notlike "*])*") {
    throw "Unexpected paragraph content at insertion anchor: $($insertBefore.Range.Text)"
}

$insertRange = $insertBefore.Range
$insertRange.Collapse(1)

$body = '<w:p><w:r><w:t>])</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>6-5 Explore $group with $unwind aggregation stage</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:lastRenderedPageBreak/><w:t>db.test.aggregate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>([</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">    // stage-1</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>{ $</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>unwind: "$friends" },</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">    // stage-2</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    {</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        $group: </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>{ _</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>id: "$friends", count: { $sum: 1 } }</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    }</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>])</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>db.test.aggregate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>([</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">    // stage-1</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>{ $</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>unwind: "$interests" },</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    // stage-2</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>{ $</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">group: { _id: "$age", </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>interestsPerAge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: {$push: "$interests"} } }</w:t></w:r></w:p><w:p/>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($xml)

Write-Output "Edit complete. Heading paragraph index=$headingIndex, insertion anchor index=$($anchorIndex + 3)"
